# Generate Report for Handoff
# Updates the Priority and Latest Handoff Datetime for the rows that were
# just handed off (091878de..., 247663df..., 2b57fcd6..., f1a943af...) on
# both the zh-cn and de-de localization-status worksheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 move from "low" priority (pending) to "ht" (handed off),
# and record the new handoff timestamp.
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-12 12:33:42"
}

# de-de: same rows, same priority flip, with the de-de handoff timestamp.
# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff datetime, so update it to the same value.
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-12 12:33:49"
    $overview.Cells.Item($r, 7).Value = "2016-08-12 12:33:49"
}
